$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy formatting (header style, currency style, total-row style) from
# column F into the new column G for every used row (1 header + 27 data + 1 total).
$ws.Range("F1:F29").Copy()
$ws.Range("G1:G29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Give the new "PRESUPUESTO" column the same width (17 characters) as in the
# target workbook.
$ws.Columns.Item(7).ColumnWidth = 16.17

# Header label
$ws.Range("G1").Value = "PRESUPUESTO"

# Data rows: budget value defaults to 0 for every advisor/client row
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

# Totals row
$ws.Range("G29").Value = 0
